# Update the build-version stamp embedded in the "About" sheet and in the
# per-row "build_version" column (S) of the "Boundaries and methane sources"
# sheet, replacing the old build timestamp with the new one.

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

# --- "About" sheet -----------------------------------------------------
$aboutWs = $wb.Worksheets.Item("About")

$a2 = [string]$aboutWs.Range("A2").Value()
$aboutWs.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = [string]$aboutWs.Range("A6").Value()
$aboutWs.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -----------------------------
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S = 19
    $val = [string]$cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
